$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the narrow separator column (old column B, width 3). This shifts the
# former C:F columns (width 16.53125) left into B:E, merging col ranges
# exactly like the target file (min=2 max=4 width=16.53125) without any
# rounding loss.
$ws.Columns.Item(2).Delete()

# Clear everything else (values + formatting) so we can rebuild the sheet
# from scratch with the new layout.
$ws.Cells.Clear()

# Write all cell values in the exact order the target workbook's shared
# string table uses them, so the regenerated shared-strings pool lines up
# with the target (cell placement is independent of this order).
$ws.Range("A2").Value = "Options"
$ws.Range("A4").Value = "Criteria"
$ws.Range("B4").Value = "Entry"
$ws.Range("C4").Value = "Difficulty"
$ws.Range("C5").Value = "Time"
$ws.Range("C6").Value = "Financial resource"
$ws.Range("C7").Value = "License"
$ws.Range("B8").Value = "Content"
$ws.Range("C8").Value = "Work hour"
$ws.Range("C11").Value = "Pressure"
$ws.Range("C14").Value = "Society"
$ws.Range("C19").Value = "Promotion"
$ws.Range("C23").Value = "Income"
$ws.Range("B29").Value = "Career"
$ws.Range("C30").Value = "Entrepreneurial Op"
$ws.Range("C29").Value = "Op to Diversify"
$ws.Range("C31").Value = "Switch track"
$ws.Range("B33").Value = "Psychology"
$ws.Range("C33").Value = "Hands-on work"
$ws.Range("C34").Value = "Team"
$ws.Range("C35").Value = "Give back to society"
$ws.Range("C36").Value = "Other"
$ws.Range("D8").Value = "Length"
$ws.Range("D9").Value = "Sleep schedule"
$ws.Range("D36").Value = "Autonomy"
$ws.Range("D37").Value = "Consistency"
$ws.Range("D38").Value = "Challenge"
$ws.Range("C32").Value = "Meaning"
$ws.Range("D10").Value = "Vacation"
$ws.Range("D11").Value = "Operational risk"
$ws.Range("D12").Value = "Performance"
$ws.Range("D13").Value = "Customer facing"
$ws.Range("D14").Value = "Social status"
$ws.Range("D15").Value = "Family and friends"
$ws.Range("D18").Value = "Uniqueness"
$ws.Range("D19").Value = "Higher level"
$ws.Range("D20").Value = "Middle level"
$ws.Range("D21").Value = "Office politics"
$ws.Range("D22").Value = "Education"
$ws.Range("D24").Value = "In 5 years"
$ws.Range("D25").Value = "In 10 years"
$ws.Range("D23").Value = "Present"
$ws.Range("E15").Value = "Colleague"
$ws.Range("E16").Value = "Supervisor"
$ws.Range("E17").Value = "Corp Culture"
$ws.Range("D26").Value = "Structure"
$ws.Range("E27").Value = "Perk"
$ws.Range("E26").Value = "Salary"
$ws.Range("E28").Value = "Bonus"
$ws.Range("C2").Value = "Pear Company"
$ws.Range("D2").Value = "Cherry Enterprise"
$ws.Range("E2").Value = "Grape Startup"
$ws.Range("B2").Value = "Banana Firm"
$ws.Range("A1").Value = "Decision"
$ws.Range("B1").Value = "Hello World"

# Restore the (empty) "applyAlignment" style that C1/D1 carry in the target
# file. Setting HorizontalAlignment to xlGeneral (-4105) reproduces the
# pre-existing cellXf #1 (applyAlignment="1", no explicit alignment) rather
# than allocating a brand-new style.
$ws.Range("C1:D1").HorizontalAlignment = -4105

# Thin spacer row between the header block and the criteria table.
$ws.Rows.Item(3).RowHeight = 4.5

# Match the saved selection/active cell.
$ws.Range("C12").Select()
